# Applies the LOB1274 content reshuffle described by the diff.
#
# The edit does not add/remove paragraphs; it only moves text between
# paragraph "slots" (each slot keeps its own paragraph/run formatting –
# pStyle, bold labels, italics, etc.). We therefore:
#   1. Capture every piece of text that needs to move, from its ORIGINAL
#      location, before making any changes.
#   2. Write the captured text into its new location.
#
# Paragraph indices (1-based, stable across the edits below since the
# total paragraph count never changes):
#   6  - Objetivos (PT)
#   7  - Objetivos (EN, italic)
#   9  - Docente(s) answer (ListBullet)
#   11 - Programa resumido (PT, with line break)
#   12 - Programa resumido (EN, italic, with line break)
#   14 - Programa (detailed, plain)
#   15 - Programa (detailed, italic) -- unchanged by this edit
#   17 - Avaliação answers (Método / Critério / Norma de recuperação)
#   19 - Bibliografia content (big block)

$d = $word.ActiveDocument

# Paragraph.Range.Text includes the trailing paragraph-mark (chr 13);
# strip it before re-using the string elsewhere, otherwise assigning it
# into another paragraph's Range.Text inserts an extra paragraph break.
function Trim-Para([string]$s) {
    return $s.TrimEnd([char]13)
}

# ---- 1. Capture original text from every source slot -----------------
$srcP6  = Trim-Para $d.Paragraphs.Item(6).Range.Text
$srcP7  = Trim-Para $d.Paragraphs.Item(7).Range.Text
$srcP9  = Trim-Para $d.Paragraphs.Item(9).Range.Text
$srcP11 = Trim-Para $d.Paragraphs.Item(11).Range.Text
$srcP12 = Trim-Para $d.Paragraphs.Item(12).Range.Text
$srcP14 = Trim-Para $d.Paragraphs.Item(14).Range.Text
$srcP19 = Trim-Para $d.Paragraphs.Item(19).Range.Text

$srcMetodoAns   = "Aulas teóricas expositivas e atividades em grupo."
$srcCriterioAns = "Média ponderada de provas e atividades."
$srcNormaAns    = "1 (uma) prova escrita"

# ---- 2. Write captured text into the new slots ------------------------

# Objetivos (PT) slot now holds the old "Programa resumido" (PT) text.
$d.Paragraphs.Item(6).Range.Text = $srcP11

# Objetivos (EN) slot now holds the old "Programa resumido" (EN) text.
$d.Paragraphs.Item(7).Range.Text = $srcP12

# Docente(s) answer slot now holds the old Objetivos (PT) text.
$d.Paragraphs.Item(9).Range.Text = $srcP6

# Programa resumido (PT) slot now holds the old detailed "Programa" (plain) text.
$d.Paragraphs.Item(11).Range.Text = $srcP14

# Programa resumido (EN) slot now holds the old Objetivos (EN) text.
$d.Paragraphs.Item(12).Range.Text = $srcP7

# Programa (plain) slot now holds the old "Método:" answer text.
$d.Paragraphs.Item(14).Range.Text = $srcMetodoAns

# Avaliação answers shift by one: Método <- Critério, Critério <- Norma,
# Norma <- old Bibliografia content block. Work within Paragraph 17's own
# range so the bold labels ("Método: ", "Critério: ", "Norma de
# recuperação: ") are left untouched.
$p17 = $d.Paragraphs.Item(17)
$r = $d.Range($p17.Range.Start, $p17.Range.End)
$r.Find.Execute($srcNormaAns, $true, $false, $false, $false, $false, $true, 1, $false, $srcP19, 2) | Out-Null

$r = $d.Range($p17.Range.Start, $p17.Range.End)
$r.Find.Execute($srcCriterioAns, $true, $false, $false, $false, $false, $true, 1, $false, $srcNormaAns, 2) | Out-Null

$r = $d.Range($p17.Range.Start, $p17.Range.End)
$r.Find.Execute($srcMetodoAns, $true, $false, $false, $false, $false, $true, 1, $false, $srcCriterioAns, 2) | Out-Null

# Bibliografia content slot now holds the old Docente(s) answer text.
$d.Paragraphs.Item(19).Range.Text = $srcP9
